$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 7806
$ws.Range("E2").Value = 509
$ws.Range("F2").Value = 509
$ws.Range("G2").Value = -535
$ws.Range("H2").Value = -522
$ws.Range("I2").Value = -381
$ws.Range("J2").Value = -141
$ws.Range("K2").Value = 29779
$ws.Range("L2").Value = 25317
$ws.Range("M2").Value = 4463
$ws.Range("N2").Value = 3665
$ws.Range("O2").Value = 798
$ws.Range("P2").Value = 2181
$ws.Range("Q2").Value = 582
$ws.Range("R2").Value = -1444
$ws.Range("S2").Value = 966
$ws.Range("T2").Value = 938
$ws.Range("U2").Value = -356
$ws.Range("V2").Value = 3485
$ws.Range("W2").Value = 6.52
$ws.Range("X2").Value = -6.68
$ws.Range("Y2").Value = -6.89
$ws.Range("Z2").Value = -1.76
$ws.Range("AA2").Value = 567.3
$ws.Range("AB2").Value = 222.2
$ws.Range("AC2").Value = -858
$ws.Range("AD2").Value = -12.43
$ws.Range("AE2").Value = 8889
$ws.Range("AF2").Value = 1.2
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 44399670
# Row 3
$ws.Range("D3").Value = 5540
$ws.Range("E3").Value = 435
$ws.Range("F3").Value = 595
$ws.Range("G3").Value = -84
$ws.Range("H3").Value = -106
$ws.Range("I3").Value = -119
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 19741
$ws.Range("L3").Value = 14056
$ws.Range("M3").Value = 5686
$ws.Range("N3").Value = 4751
$ws.Range("O3").Value = 935
$ws.Range("P3").Value = 1930
$ws.Range("Q3").Value = 1112
$ws.Range("R3").Value = -14
$ws.Range("S3").Value = -859
$ws.Range("T3").Value = 307
$ws.Range("U3").Value = 805
$ws.Range("V3").Value = 4803
$ws.Range("W3").Value = 7.85
$ws.Range("X3").Value = -1.92
$ws.Range("Y3").Value = -2.82
$ws.Range("Z3").Value = -0.43
$ws.Range("AA3").Value = 247.2
$ws.Range("AB3").Value = 217.69
$ws.Range("AC3").Value = -354
$ws.Range("AD3").Value = -19.22
$ws.Range("AE3").Value = 13531
$ws.Range("AF3").Value = 0.5
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 39282593
# Row 4
$ws.Range("D4").Value = 7724
$ws.Range("E4").Value = 401
$ws.Range("F4").Value = 623
$ws.Range("G4").Value = 159
$ws.Range("H4").Value = -69
$ws.Range("I4").Value = -33
$ws.Range("J4").Value = -36
$ws.Range("K4").Value = 20067
$ws.Range("L4").Value = 13884
$ws.Range("M4").Value = 6184
$ws.Range("N4").Value = 5216
$ws.Range("O4").Value = 968
$ws.Range("P4").Value = 2318
$ws.Range("Q4").Value = 183
$ws.Range("R4").Value = -702
$ws.Range("S4").Value = 509
$ws.Range("T4").Value = 440
$ws.Range("U4").Value = -257
$ws.Range("V4").Value = 2183
$ws.Range("W4").Value = 5.19
$ws.Range("X4").Value = -0.9
$ws.Range("Y4").Value = -0.67
$ws.Range("Z4").Value = -0.35
$ws.Range("AA4").Value = 224.53
$ws.Range("AB4").Value = 187.53
$ws.Range("AC4").Value = -72
$ws.Range("AD4").Value = -91.11
$ws.Range("AE4").Value = 12125
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 47183679
# Row 5
$ws.Range("D5").Value = 6868
$ws.Range("E5").Value = 256
$ws.Range("F5").Value = 256
$ws.Range("G5").Value = -631
$ws.Range("H5").Value = -867
$ws.Range("I5").Value = -549
$ws.Range("J5").Value = -318
$ws.Range("K5").Value = 16294
$ws.Range("L5").Value = 11063
$ws.Range("M5").Value = 5231
$ws.Range("N5").Value = 4670
$ws.Range("O5").Value = 561
$ws.Range("P5").Value = 2318
$ws.Range("Q5").Value = 567
$ws.Range("R5").Value = -310
$ws.Range("S5").Value = -248
$ws.Range("T5").Value = 256
$ws.Range("U5").Value = 311
$ws.Range("V5").Value = 2068
$ws.Range("W5").Value = 3.73
$ws.Range("X5").Value = -12.63
$ws.Range("Y5").Value = -11.11
$ws.Range("Z5").Value = -4.77
$ws.Range("AA5").Value = 211.48
$ws.Range("AB5").Value = 164.75
$ws.Range("AC5").Value = -1164
$ws.Range("AD5").Value = -3.86
$ws.Range("AE5").Value = 10856
$ws.Range("AF5").Value = 0.41
$ws.Range("AG5").Value = 49
$ws.Range("AH5").Value = 1.09
$ws.Range("AI5").Value = -3.84
$ws.Range("AJ5").Value = 47183679
# Row 6
$ws.Range("D6").Value = 6987
$ws.Range("E6").Value = 381
$ws.Range("F6").Value = 381
$ws.Range("G6").Value = -159
$ws.Range("H6").Value = -415
$ws.Range("I6").Value = -390
$ws.Range("K6").Value = 13474
$ws.Range("L6").Value = 9154
$ws.Range("M6").Value = 4320
$ws.Range("N6").Value = 4090
$ws.Range("P6").Value = 2318
$ws.Range("Q6").Value = -38
$ws.Range("R6").Value = 504
$ws.Range("S6").Value = -359
$ws.Range("T6").Value = 141
$ws.Range("U6").Value = -179
$ws.Range("V6").Value = 1143
$ws.Range("W6").Value = 5.46
$ws.Range("X6").Value = -5.94
$ws.Range("Y6").Value = -8.9
$ws.Range("Z6").Value = -2.79
$ws.Range("AA6").Value = 211.91
$ws.Range("AB6").Value = 142.92
$ws.Range("AC6").Value = -826
$ws.Range("AD6").Value = -5.95
$ws.Range("AE6").Value = 9735
$ws.Range("AF6").Value = 0.5
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 47183679
# Row 7
$ws.Range("D7").Value = 5430
$ws.Range("E7").Value = 190
$ws.Range("G7").Value = 410
$ws.Range("H7").Value = 460
$ws.Range("I7").Value = 450
$ws.Range("K7").Value = 6900
$ws.Range("L7").Value = 2140
$ws.Range("M7").Value = 4770
$ws.Range("N7").Value = 4540
$ws.Range("P7").Value = 2360
$ws.Range("Q7").Value = 430
$ws.Range("R7").Value = -380
$ws.Range("S7").Value = -140
$ws.Range("T7").Value = 130
$ws.Range("W7").Value = 3.5
$ws.Range("X7").Value = 8.470000000000001
$ws.Range("Y7").Value = 10.43
$ws.Range("Z7").Value = 4.52
$ws.Range("AA7").Value = 44.86
$ws.Range("AC7").Value = 954
$ws.Range("AD7").Value = 3.67
$ws.Range("AE7").Value = 10807
$ws.Range("AF7").Value = 0.32
$ws.Range("AG7").Value = 100
$ws.Range("AH7").Value = 2.86
$ws.Range("AI7").Value = 10.48
# Row 8
$ws.Range("D8").Value = 5400
$ws.Range("E8").Value = 170
$ws.Range("G8").Value = 350
$ws.Range("H8").Value = 280
$ws.Range("I8").Value = 260
$ws.Range("K8").Value = 7090
$ws.Range("L8").Value = 2090
$ws.Range("M8").Value = 5000
$ws.Range("N8").Value = 4760
$ws.Range("P8").Value = 2360
$ws.Range("Q8").Value = 160
$ws.Range("R8").Value = -140
$ws.Range("S8").Value = 80
$ws.Range("T8").Value = 140
$ws.Range("W8").Value = 3.15
$ws.Range("X8").Value = 5.18
$ws.Range("Y8").Value = 5.59
$ws.Range("Z8").Value = 4
$ws.Range("AA8").Value = 41.8
$ws.Range("AC8").Value = 551
$ws.Range("AD8").Value = 6.35
$ws.Range("AE8").Value = 11331
$ws.Range("AF8").Value = 0.31
$ws.Range("AG8").Value = 100
$ws.Range("AH8").Value = 2.86
$ws.Range("AI8").Value = 18.15
# Row 9
$ws.Range("D9").Value = 5330
$ws.Range("E9").Value = 190
$ws.Range("G9").Value = 380
$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 290
$ws.Range("K9").Value = 7280
$ws.Range("L9").Value = 2020
$ws.Range("M9").Value = 5260
$ws.Range("N9").Value = 5010
$ws.Range("P9").Value = 2360
$ws.Range("Q9").Value = 180
$ws.Range("R9").Value = -150
$ws.Range("S9").Value = 90
$ws.Range("T9").Value = 150
$ws.Range("W9").Value = 3.56
$ws.Range("X9").Value = 5.63
$ws.Range("Y9").Value = 5.94
$ws.Range("Z9").Value = 5.69
$ws.Range("AA9").Value = 38.4
$ws.Range("AC9").Value = 615
$ws.Range("AD9").Value = 5.69
$ws.Range("AE9").Value = 11926
$ws.Range("AF9").Value = 0.29
$ws.Range("AG9").Value = 100
$ws.Range("AH9").Value = 2.86
$ws.Range("AI9").Value = 16.27

# Remove cells for non-controlling interest in estimate years (merged into T column)
$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()